$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: update/clear individual cell values (row numbers as in the original sheet) ---

# E3 was empty -> now has a value
$ws.Range("E3").Value = -5.7

# F4 had a value -> now cleared (missing)
$ws.Range("F4").ClearContents()

# E5 had a value -> now cleared (missing)
$ws.Range("E5").ClearContents()

# F9, F10, F11, F12 were empty -> now have values
$ws.Range("F9").Value = 17.26
$ws.Range("F10").Value = 16.43
$ws.Range("F11").Value = 17.65
$ws.Range("F12").Value = 17.45

# F15, F17, F18, F20 had values -> now cleared (missing)
$ws.Range("F15").ClearContents()
$ws.Range("F17").ClearContents()
$ws.Range("F18").ClearContents()
$ws.Range("F20").ClearContents()

# E21 was empty -> now has a value
$ws.Range("E21").Value = -8.699999999999999

# E23 had a value -> now cleared (missing)
$ws.Range("E23").ClearContents()

# --- Step 2: remove the rows for "RM 232" (row 26) and "SC 92" (row 28) entirely ---
# Delete the higher-numbered row first so the other row index stays valid.
$ws.Rows(28).Delete()
$ws.Rows(26).Delete()

# --- Step 3: after the rows above shift up, a couple of formerly-empty cells in the
#     shifted "SC 132" / "SC 193" rows (now rows 31 and 32) gain values ---
$ws.Range("F31").Value = 17.18
$ws.Range("E32").Value = -6.4
$ws.Range("F32").Value = 17.39

Write-Output "edit complete"
